$d = $word.ActiveDocument

function ReplaceOnce($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Not found: $findText"
    }
}

# 1. "Het netwerk policy geld" -> "geldt" (Scope under "Netwerk")
ReplaceOnce "Het netwerk policy geld voor zowel persoonsleden, stude" "Het netwerk policy geldt voor zowel persoonsleden, stude"

# 2. "porno" -> "pornografisch materiaal"
ReplaceOnce "illegaal downloaden en porno zijn ten striktste verboden." "illegaal downloaden en pornografisch materiaal zijn ten striktste verboden."

# 3. "zullen hard tegen worden opgetreden." -> "zijn verboden."
ReplaceOnce "Programma’s en of tools die het netwerk in gevaar kunnen brengen zullen hard tegen worden opgetreden." "Programma’s en of tools die het netwerk in gevaar kunnen brengen zijn verboden."

# 4. "een gemakkelijke manier ... kan maken." -> "een gemakkelijke en vlotte manier ... kunnen maken."
ReplaceOnce " een gemakkelijke manier hiervan gebruik kan maken." " een gemakkelijke en vlotte manier hiervan gebruik kunnen maken."

# 5. "De printer policy geld" -> "geldt"
ReplaceOnce "De printer policy geld op dezelfde" "De printer policy geldt op dezelfde"

# 6. "en/of zo de documenten" -> "en/of de documenten"
ReplaceOnce "en/of zo de documenten zo optimaal" "en/of de documenten zo optimaal"

# 7. "haal direct iemand" -> "haal iemand"
ReplaceOnce "haal direct iemand die weet" "haal iemand die weet"

# 8. "Het netwerk policy geld voor zowel persoonsleden als studenten" -> "geldt"
ReplaceOnce "Het netwerk policy geld voor zowel persoonsleden als studenten" "Het netwerk policy geldt voor zowel persoonsleden als studenten"

# 9. Remove the _GoBack bookmark that currently sits between "Ieder personeelslid" and
#    "/student ...". A plain Find&Replace over that span rewrites the run(s) in place and
#    (as a side effect, same as real Word re-typing text there) drops the now-orphaned
#    bookmark markers without disturbing anything else.
ReplaceOnce "Ieder personeelslid/student die een gebruikersnaam en wachtwoord heeft." "Ieder personeelslid/student die een gebruikersnaam en wachtwoord heeft."

# 10. "Kies een veilig wachtwoord. Minstens bestaande uit hoofdletters, cijfers en symbolen." ->
#     "...met een lengte van tenminste 8 karakters." (and the _GoBack bookmark re-appears here,
#     between "wacht" and "woord", matching where the author was last editing.)
ReplaceOnce "Kies een veilig wachtwoord. Minstens bestaande uit hoofdletters, cijfers en symbolen." "Kies een veilig wachtwoord. Minstens bestaande uit hoofdletters, cijfers en symbolen met een lengte van tenminste 8 karakters."

$r3 = $d.Content
$found3 = $r3.Find.Execute("Kies een veilig wacht", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $pos = $r3.End
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Host "Done"
